$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '30.241.15'
Set-TextValue $ws.Range('E2') '  -0.23%  '
Set-TextValue $ws.Range('D3') '1.926.34'
Set-TextValue $ws.Range('E4') '  +0.05%  '
Set-TextValue $ws.Range('D5') '248.66'
Set-TextValue $ws.Range('E5') '  +0.02%  '
Set-TextValue $ws.Range('D6') '0.7126'
Set-TextValue $ws.Range('E6') '  -1.15%  '
Set-TextValue $ws.Range('E7') '  +0.06%  '
Set-TextValue $ws.Range('E8') '  -2.96%  '
Set-TextValue $ws.Range('D9') '27.43'
Set-TextValue $ws.Range('E9') '  -1.60%  '
Set-TextValue $ws.Range('D10') '0.07060'
Set-TextValue $ws.Range('E10') '  +1.94%  '
Set-TextValue $ws.Range('D11') '0.7914'
Set-TextValue $ws.Range('E11') '  -1.37%  '
Set-TextValue $ws.Range('D12') '0.07959'
Set-TextValue $ws.Range('E12') '  -1.33%  '
Set-TextValue $ws.Range('D13') '1.930.96'
Set-TextValue $ws.Range('E13') '  +0.07%  '
Set-TextValue $ws.Range('D14') '5.377'
Set-TextValue $ws.Range('E14') '  -0.61%  '
Set-TextValue $ws.Range('D15') '94.77'
Set-TextValue $ws.Range('E15') '  +0.06%  '
Set-TextValue $ws.Range('D16') '14.61'
Set-TextValue $ws.Range('E16') '  +0.78%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D17') '30.262.97'
Set-TextValue $ws.Range('E17') '  -0.12%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D18') '258.04'
Set-TextValue $ws.Range('E18') '  +2.03%  '
Set-TextValue $ws.Range('D19') '0.000008044'
Set-TextValue $ws.Range('E19') '  -3.73%  '
Set-TextValue $ws.Range('D20') '5.754'
Set-TextValue $ws.Range('E20') '  -0.84%  '
Set-TextValue $ws.Range('D21') '2.184.95'
Set-TextValue $ws.Range('E21') '  -0.01%  '
Set-TextValue $ws.Range('E22') '  +0.04%  '
Set-TextValue $ws.Range('E23') '  +0.03%  '
Set-TextValue $ws.Range('D24') '6.850'
Set-TextValue $ws.Range('E24') '  -0.20%  '
Set-TextValue $ws.Range('D25') '9.525'
Set-TextValue $ws.Range('E25') '  -1.99%  '
Set-TextValue $ws.Range('D26') '165.97'
Set-TextValue $ws.Range('E26') '  +4.13%  '
Set-TextValue $ws.Range('D27') '19.08'
Set-TextValue $ws.Range('E27') '  -0.32%  '
Set-TextValue $ws.Range('D28') '2.259'
Set-TextValue $ws.Range('E28') '  -5.95%  '
Set-TextValue $ws.Range('D29') '0.1257'
Set-TextValue $ws.Range('E29') '  -5.79%  '
Set-TextValue $ws.Range('E30') '  +1.25%  '
Set-TextValue $ws.Range('E31') '  -1.67%  '
Set-TextValue $ws.Range('D32') '4.387'
Set-TextValue $ws.Range('E32') '  -0.37%  '
Set-TextValue $ws.Range('D33') '4.119'
Set-TextValue $ws.Range('E33') '  -1.68%  '
Set-TextValue $ws.Range('D34') '0.05131'
Set-TextValue $ws.Range('E34') '  +0.47%  '
Set-TextValue $ws.Range('E35') '  +3.98%  '
Set-TextValue $ws.Range('D36') '0.7436'
Set-TextValue $ws.Range('E36') '  +0.52%  '
Set-TextValue $ws.Range('D37') '2.766'
Set-TextValue $ws.Range('E37') '  +0.89%  '
Set-TextValue $ws.Range('D38') '0.01957'
Set-TextValue $ws.Range('E38') '  -0.50%  '
Set-TextValue $ws.Range('D39') '2.796'
Set-TextValue $ws.Range('E39') '  -1.17%  '
Set-TextValue $ws.Range('D40') '77.36'
Set-TextValue $ws.Range('E40') '  -1.88%  '
Set-TextValue $ws.Range('D41') '6.355'
Set-TextValue $ws.Range('E41') '  -3.54%  '
Set-TextValue $ws.Range('D42') '0.4493'
Set-TextValue $ws.Range('E42') '  +0.71%  '
Set-TextValue $ws.Range('D43') '1.987'
Set-TextValue $ws.Range('E43') '  -0.15%  '
Set-TextValue $ws.Range('D44') '0.8451'
Set-TextValue $ws.Range('E44') '  +1.09%  '
Set-TextValue $ws.Range('D45') '1.000'
Set-TextValue $ws.Range('E45') '  -0.03%  '
Set-TextValue $ws.Range('D46') '100.65'
Set-TextValue $ws.Range('E46') '  -1.36%  '
Set-TextValue $ws.Range('D47') '9.740'
Set-TextValue $ws.Range('E47') '  -0.55%  '
Set-TextValue $ws.Range('D48') '7.415'
Set-TextValue $ws.Range('E48') '  +1.65%  '
Set-TextValue $ws.Range('D49') '36.54'
Set-TextValue $ws.Range('E49') '  +0.08%  '
Set-TextValue $ws.Range('D50') '0.06108'
Set-TextValue $ws.Range('E50') '  +2.58%  '
Set-TextValue $ws.Range('D51') '0.4187'
Set-TextValue $ws.Range('E51') '  +2.62%  '
